# BookDatabase.xlsx edit
# Commit: "Changed from Google database to wcat. Changed to store values as
#          numbers over text."
#
# Sheet "Book Inventory" (sheet1): the inventory previously held three
# duplicate rows for "The Hunger Games". We now keep only two book rows:
#   - Row 2 becomes "Where The Wild Things Are" by Maurice Sendak (wcat)
#   - Row 3 keeps "The Hunger Games" by Suzanne Collins
#   - The extra duplicate row 4 is removed
#   - Both rows get numeric Total Quantity / In Stock values of 1
#
# Sheet "Check Out-In" (sheet2) is unaffected content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Book Inventory")

$ws.Activate()

# Drop the extra duplicate row (the inventory only has 2 distinct books now).
$ws.Rows.Item(4).Delete()

# Row 2: replace the old "Hunger Games" placeholder with the new book.
$ws.Range("A2").Value = "Where The Wild Things Are"
$ws.Range("B2").Value = "Maurice Sendak"

# ISBNs are stored as text so leading zeros survive; apply a Text format
# before writing the value, then drop back to the Normal cell style so no
# visible formatting difference remains on the cell.
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0064431789"
$ws.Range("C2").Style = "Normal"

# Quantities are now stored as real numbers instead of text.
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3 keeps "The Hunger Games" / "Suzanne Collins" / "0439023483" (already
# present as shared text), just add the numeric quantity columns.
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 1

# Reflect the author's last selection before saving.
[void]$ws.Range("A2").Select()
